$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1400
$ws.Range("J9").Value = 499
$ws.Range("L9").Value = 499
$ws.Range("N9").Value = -837
$ws.Range("H17").Value = 1367.4615
$ws.Range("J17").Value = 1367.4615
$ws.Range("L17").Value = 4102.3845
$ws.Range("N17").Value = -4438.3845
$ws.Range("H53").Value = 592.4
$ws.Range("J53").Value = 721.75
$ws.Range("L53").Value = 721.75
$ws.Range("N53").Value = -1995.75
$ws.Range("H98").Value = 497.5
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("H122").Value = 497.5
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = $null
$ws.Range("H135").Value = 1490.4166
$ws.Range("I135").Value = 1490.4166
$ws.Range("K135").Value = 13413.7494
$ws.Range("M135").Value = -10878.7494
$ws.Range("H137").Value = 2200.5557
$ws.Range("I137").Value = 1097.1666
$ws.Range("J137").Value = 4407.3335
$ws.Range("K137").Value = 3291.4998
$ws.Range("L137").Value = 13222.0005
$ws.Range("M137").Value = -741.4998000000001
$ws.Range("N137").Value = -18322.0005
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2546.0312
$ws.Range("I32").Value = 2464.6206
$ws.Range("K32").Value = 2464.6206
$ws.Range("M32").Value = -2177.6206
$ws.Range("H61").Value = 3483.8572
$ws.Range("I61").Value = 3905.5625
$ws.Range("K61").Value = 3905.5625
$ws.Range("M61").Value = -3693.5625
$ws.Range("H122").Value = 1758.6
$ws.Range("J122").Value = 1198
$ws.Range("L122").Value = 3594
$ws.Range("N122").Value = -8494
$ws.Range("H132").Value = 1956.0278
$ws.Range("I132").Value = 1755.875
$ws.Range("J132").Value = 3557.25
$ws.Range("K132").Value = 5267.625
$ws.Range("L132").Value = 10671.75
$ws.Range("M132").Value = -2737.625
$ws.Range("N132").Value = -15731.75
$ws.Range("H136").Value = 3483.8572
$ws.Range("I136").Value = 3905.5625
$ws.Range("K136").Value = 11716.6875
$ws.Range("M136").Value = -9166.6875
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 598.8
$ws.Range("I22").Value = 573.5
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 573.5
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = -400.5
$ws.Range("N22").Value = -1046
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2960.75
$ws.Range("J94").Value = 2956
$ws.Range("L94").Value = 2956
$ws.Range("N94").Value = -3858
$ws.Range("H99").Value = 3899.75
$ws.Range("I99").Value = 4000
$ws.Range("J99").Value = 3799.5
$ws.Range("K99").Value = 4000
$ws.Range("L99").Value = 3799.5
$ws.Range("M99").Value = -2502
$ws.Range("N99").Value = -6795.5
$ws.Range("H126").Value = 3899.75
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 3799.5
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 11398.5
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -16338.5
$ws.Range("H132").Value = 3761.2727
$ws.Range("I132").Value = 2483.4285
$ws.Range("K132").Value = 7450.2855
$ws.Range("M132").Value = -4920.2855
$ws.Range("H134").Value = 4091.4614
$ws.Range("I134").Value = 4091.4614
$ws.Range("K134").Value = 12274.3842
$ws.Range("M134").Value = -9739.3842
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 325.7143
$ws.Range("I12").Value = 307.75
$ws.Range("J12").Value = 349.66666
$ws.Range("K12").Value = 923.25
$ws.Range("L12").Value = 1048.99998
$ws.Range("M12").Value = -750.25
$ws.Range("N12").Value = -1394.99998
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3069.8572
$ws.Range("I126").Value = 3069.8572
$ws.Range("K126").Value = 9209.571599999999
$ws.Range("M126").Value = -6739.571599999999
$ws.Range("H132").Value = 2959.6667
$ws.Range("I132").Value = 2480.4546
$ws.Range("J132").Value = 4277.5
$ws.Range("K132").Value = 7441.3638
$ws.Range("L132").Value = 12832.5
$ws.Range("M132").Value = -4911.3638
$ws.Range("N132").Value = -17892.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7499.5
$ws.Range("I7").Value = 7499.5
$ws.Range("K7").Value = 7499.5
$ws.Range("M7").Value = -7387.5
$ws.Range("H40").Value = 2600
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = $null
$ws.Range("H82").Value = 1488
$ws.Range("I82").Value = 1799.6666
$ws.Range("J82").Value = 1176.3334
$ws.Range("K82").Value = 1799.6666
$ws.Range("L82").Value = 1176.3334
$ws.Range("M82").Value = -1438.6666
$ws.Range("N82").Value = -1898.3334
$ws.Range("H85").Value = 1488
$ws.Range("I85").Value = 1799.6666
$ws.Range("J85").Value = 1176.3334
$ws.Range("K85").Value = 1799.6666
$ws.Range("L85").Value = 1176.3334
$ws.Range("M85").Value = -551.6666
$ws.Range("N85").Value = -3672.3334
$ws.Range("H88").Value = 19999
$ws.Range("J88").Value = 19999
$ws.Range("L88").Value = 19999
$ws.Range("N88").Value = -20855
$ws.Range("H91").Value = 19999
$ws.Range("J91").Value = 19999
$ws.Range("L91").Value = 19999
$ws.Range("N91").Value = -22963
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("H126").Value = 7499.5
$ws.Range("I126").Value = 7499.5
$ws.Range("K126").Value = 22498.5
$ws.Range("M126").Value = -20028.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7499.5
$ws.Range("I62").Value = 9999
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 9999
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -9375
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 7499.5
$ws.Range("I65").Value = 9999
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 49995
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -46875
$ws.Range("N65").Value = -31240
$ws.Range("H100").Value = 819.8
$ws.Range("I100").Value = 819.8
$ws.Range("K100").Value = 1639.6
$ws.Range("M100").Value = -1098.6
$ws.Range("H126").Value = 2078.25
$ws.Range("I126").Value = 152
$ws.Range("K126").Value = 456
$ws.Range("M126").Value = 2014
